$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Body text (Simplified Chinese -> Traditional Chinese) ---

Replace-Text "英语" "英語"
Replace-Text " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" " / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語"
Replace-Text "简要" "簡介"
Replace-Text "发给在目标国家已确认出席但未在截止日期前提交文件的合作伙伴的邮件。 我们将撤回他们的邀请。 将通过 customer.io 发送" "發送給在目標國家的合作夥伴的電子郵件，這些合作夥伴已回應參加，但在截止日期前未提交文件。 我們將取消他們的邀請。 將通過 customer.io 發送"
Replace-Text "目标受众" "目標受眾"
Replace-Text "未及时提交文件的被邀请合作伙伴" "未按時提交文件的被邀請合作夥伴"
Replace-Text "主题行" "主題行"
Replace-Text "[活动名称]" "[活動名稱]"
Replace-Text " 注册" " 註冊"
Replace-Text "没有及时收到您的文件" "沒有及時收到您的文件"
Replace-Text "[PARTNER NAME]" "[合作夥伴姓名]"
Replace-Text "We didn’t receive your documents by the deadline (" "截止日期（"
Replace-Text "[DD Mmm YYYY]" "[日月年]"
Replace-Text "). 很遗憾，无法为您办理 " "）前沒有收到您的文件。 很遺憾，無法為您辦理 "
Replace-Text " 的注册手续。" " 的註冊手續。"
Replace-Text "衷心祝愿您一切顺利，并希望在下一次 " "衷心祝愿您一切順利，並希望在下一次 "
Replace-Text "会议/研讨会/联盟会员旅行" "會議/研討會/聯盟會員旅行"
Replace-Text "中见到您。" "中見到您。"
Replace-Text "如有任何疑问，请通过 " "如有任何疑問，請通過 "
Replace-Text "[电子邮件地址]" "[電子郵件地址]"
Replace-Text "[WHATSAPP 号码]" "[WHATSAPP 號碼]"
Replace-Text " (WhatsApp) 联系您的区域经理 " " (WhatsApp) 聯繫您的區域經理 "

# --- Comments text (Simplified Chinese -> Traditional Chinese) ---
$c = $d.Comments
for ($i = 1; $i -le $c.Count; $i++) {
    $c.Item($i).Range.Find.Execute("选择任一", $true, $false, $false, $false, $false, $true, 1, $false, "選擇其中一個", 2) | Out-Null
}
